# NYPD CompStat weekly report refresh: new week's crime data collected.
# Updates the report header (volume/week-of-number text + date range) and
# the "This Week" statistics table (rows 14-30, 33) with the newly
# collected weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# xlPasteSpecial constants used below
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# A handful of cells in the "This Week" table flip between being a plain
# number (style 14/15, numeric <v>) and the literal text placeholders
# "0" / "***.*" (style 13, shared string). Excel's Range.Value setter
# auto-detects numeric-looking strings and stores them as numbers, so a
# straight assignment can't produce a *text* "0" in a General-formatted
# cell. Instead we copy both the value and the number format from an
# existing cell that already holds the exact text/format we need.
function Copy-CellValueAndFormat($srcRange, $dstRange) {
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial($xlPasteValues) | Out-Null
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial($xlPasteFormats) | Out-Null
}

# Turn a cell into the literal text "0" (style 13), copying from a cell
# that already has exactly that text+style.
function Set-TextZero([string]$cellRef) {
    Copy-CellValueAndFormat $ws.Range("D14") $ws.Range($cellRef)
}

# Turn a cell into the literal text "***.*" (style 13), copying from a
# cell that already has exactly that text+style.
function Set-TextStar([string]$cellRef) {
    Copy-CellValueAndFormat $ws.Range("E14") $ws.Range($cellRef)
}

# Set a cell to a plain number, first pasting in the number format from
# an existing cell with the desired style (14 = "#,##0", 15 = "#,##0.0").
function Set-NumberValue([string]$cellRef, [string]$formatSrcRef, $value) {
    $src = $ws.Range($formatSrcRef)
    $dst = $ws.Range($cellRef)
    $src.Copy() | Out-Null
    $dst.PasteSpecial($xlPasteFormats) | Out-Null
    $dst.Value = $value
}

# ---------------------------------------------------------------------
# Header: volume/number + week-covering date range
# ---------------------------------------------------------------------

$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# ---------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------
Set-TextZero "C14"
$ws.Range("N14").Value = -88.135593220339

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-NumberValue "C15" "F14" 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 34
$ws.Range("K15").Value = 9.677419354838
$ws.Range("L15").Value = 78.947368421052
$ws.Range("M15").Value = 54.545454545454
$ws.Range("N15").Value = -53.424657534246

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 370
$ws.Range("J16").Value = 347
$ws.Range("K16").Value = 6.628242074927
$ws.Range("L16").Value = 18.589743589743
$ws.Range("M16").Value = -5.612244897959
$ws.Range("N16").Value = -80.412916887241

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 19
$ws.Range("E17").Value = -5.263157894736
$ws.Range("F17").Value = 71
$ws.Range("G17").Value = 79
$ws.Range("H17").Value = -10.126582278481
$ws.Range("I17").Value = 774
$ws.Range("J17").Value = 794
$ws.Range("K17").Value = -2.518891687657
$ws.Range("L17").Value = 9.322033898305
$ws.Range("M17").Value = 108.625336927224
$ws.Range("N17").Value = -16.594827586206

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 175
$ws.Range("J18").Value = 196
$ws.Range("K18").Value = -10.714285714285
$ws.Range("L18").Value = 19.047619047619
$ws.Range("M18").Value = -2.234636871508
$ws.Range("N18").Value = -91.471734892787

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -31.578947368421
$ws.Range("I19").Value = 579
$ws.Range("J19").Value = 608
$ws.Range("K19").Value = -4.769736842105
$ws.Range("L19").Value = 17.682926829268
$ws.Range("M19").Value = 120.992366412214
$ws.Range("N19").Value = -32.674418604651

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 550
$ws.Range("F20").Value = 23
$ws.Range("H20").Value = 53.333333333333
$ws.Range("I20").Value = 229
$ws.Range("J20").Value = 190
$ws.Range("K20").Value = 20.526315789473
$ws.Range("L20").Value = -9.126984126984
$ws.Range("M20").Value = 108.181818181818
$ws.Range("N20").Value = -70.716112531969

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 11.111111111111
$ws.Range("F21").Value = 184
$ws.Range("G21").Value = 206
$ws.Range("H21").Value = -10.679611650485
$ws.Range("I21").Value = 2168
$ws.Range("J21").Value = 2191
$ws.Range("K21").Value = -1.049748973071
$ws.Range("L21").Value = 11.637487126673
$ws.Range("M21").Value = 60.830860534124
$ws.Range("N21").Value = -67.364142706608

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("D22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 26
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 18.181818181818
$ws.Range("L22").Value = -7.142857142857
$ws.Range("M22").Value = -7.142857142857

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
Set-NumberValue "D23" "F23" 2
Set-NumberValue "E23" "H22" -100
$ws.Range("F23").Value = 2
Set-NumberValue "G23" "F23" 2
Set-NumberValue "H23" "H22" 0
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = 30.769230769230
$ws.Range("L23").Value = 36
$ws.Range("M23").Value = 61.904761904761

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 61
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = -25.609756097561
$ws.Range("I24").Value = 859
$ws.Range("J24").Value = 848
$ws.Range("K24").Value = 1.297169811320
$ws.Range("L24").Value = -7.135135135135
$ws.Range("M24").Value = 3.995157384987

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -68
$ws.Range("I25").Value = 221
$ws.Range("J25").Value = 298
$ws.Range("K25").Value = -25.838926174496
$ws.Range("L25").Value = 17.553191489361

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 6.25
$ws.Range("F26").Value = 77
$ws.Range("G26").Value = 95
$ws.Range("H26").Value = -18.947368421052
$ws.Range("I26").Value = 923
$ws.Range("J26").Value = 1008
$ws.Range("K26").Value = -8.432539682539
$ws.Range("L26").Value = -13.495782567947
$ws.Range("M26").Value = -19.388646288209

# ---------------------------------------------------------------------
# Row 27 - UCR Rape
# ---------------------------------------------------------------------
Set-NumberValue "C27" "F14" 1
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 47
$ws.Range("K27").Value = 4.444444444444
$ws.Range("L27").Value = -6

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
Set-TextZero "C28"
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 82
$ws.Range("J28").Value = 85
$ws.Range("K28").Value = -3.529411764705
$ws.Range("L28").Value = 1.234567901234

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
Set-TextZero "C29"
Set-NumberValue "D29" "F14" 3
Set-NumberValue "E29" "H22" -100
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = -66.666666666666
$ws.Range("J29").Value = 62
$ws.Range("K29").Value = -59.677419354838
$ws.Range("M29").Value = -50.980392156862
$ws.Range("N29").Value = -84.848484848484

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------
Set-TextZero "C30"
Set-NumberValue "D30" "F14" 1
Set-NumberValue "E30" "H22" -100
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 47
$ws.Range("K30").Value = -55.319148936170
$ws.Range("M30").Value = -52.272727272727
$ws.Range("N30").Value = -86.092715231788

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities
# ---------------------------------------------------------------------
Set-TextZero "F33"
